$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(98, 8).Value = 2862.6155  # H98: 2664.7144 -> 2862.6155
$ws.Cells.Item(98, 9).Value = 2862.6155  # I98: 2664.7144 -> 2862.6155
$ws.Cells.Item(98, 11).Value = 2862.6155  # K98: 2664.7144 -> 2862.6155
$ws.Cells.Item(98, 13).Value = -1364.6155  # M98: -1166.7144 -> -1364.6155
$ws.Cells.Item(116, 8).Value = 2922.282  # H116: 2967.7104 -> 2922.282
$ws.Cells.Item(116, 9).Value = 2618.64  # I116: 2677.9167 -> 2618.64
$ws.Cells.Item(116, 11).Value = 2618.64  # K116: 2677.9167 -> 2618.64
$ws.Cells.Item(116, 13).Value = 823.3600000000001  # M116: 764.0832999999998 -> 823.3600000000001
$ws.Cells.Item(122, 8).Value = 2862.6155  # H122: 2664.7144 -> 2862.6155
$ws.Cells.Item(122, 9).Value = 2862.6155  # I122: 2664.7144 -> 2862.6155
$ws.Cells.Item(122, 11).Value = 8587.8465  # K122: 7994.1432 -> 8587.8465
$ws.Cells.Item(122, 13).Value = -6137.8465  # M122: -5544.1432 -> -6137.8465
$ws.Cells.Item(132, 8).Value = 11117325  # H132: 12827575 -> 11117325
$ws.Cells.Item(132, 9).Value = 12827138  # I132: 15159217 -> 12827138
$ws.Cells.Item(132, 11).Value = 38481414  # K132: 45477651 -> 38481414
$ws.Cells.Item(132, 13).Value = -38478884  # M132: -45475121 -> -38478884
$ws.Cells.Item(137, 8).Value = 1054.8572  # H137: 1054.9143 -> 1054.8572
$ws.Cells.Item(137, 9).Value = 911.8182  # I137: 911.9091 -> 911.8182
$ws.Cells.Item(137, 11).Value = 2735.4546  # K137: 2735.7273 -> 2735.4546
$ws.Cells.Item(137, 13).Value = -185.4546  # M137: -185.7273 -> -185.4546
$ws.Cells.Item(138, 8).Value = 1568.1803  # H138: 1604.5424 -> 1568.1803
$ws.Cells.Item(138, 9).Value = 1360.9286  # I138: 1417.96 -> 1360.9286
$ws.Cells.Item(138, 10).Value = 1744.0303  # J138: 1741.7354 -> 1744.0303
$ws.Cells.Item(138, 11).Value = 4082.7858  # K138: 4253.88 -> 4082.7858
$ws.Cells.Item(138, 12).Value = 5232.090899999999  # L138: 5225.206200000001 -> 5232.090899999999
$ws.Cells.Item(138, 13).Value = 1057.2142  # M138: 886.1199999999999 -> 1057.2142
$ws.Cells.Item(138, 14).Value = -15512.0909  # N138: -15505.2062 -> -15512.0909

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(61, 8).Value = 990  # H61: 945.6667 -> 990
$ws.Cells.Item(61, 9).Value = 990  # I61: 945.6667 -> 990
$ws.Cells.Item(61, 11).Value = 990  # K61: 945.6667 -> 990
$ws.Cells.Item(61, 13).Value = -778  # M61: -733.6667 -> -778
$ws.Cells.Item(74, 8).Value = 884.35895  # H74: 787.13043 -> 884.35895
$ws.Cells.Item(74, 9).Value = 502.06668  # I74: 460 -> 502.06668
$ws.Cells.Item(74, 10).Value = 2158.6667  # J74: 2341 -> 2158.6667
$ws.Cells.Item(74, 11).Value = 502.06668  # K74: 460 -> 502.06668
$ws.Cells.Item(74, 12).Value = 2158.6667  # L74: 2341 -> 2158.6667
$ws.Cells.Item(74, 13).Value = 371.93332  # M74: 414 -> 371.93332
$ws.Cells.Item(74, 14).Value = -3906.6667  # N74: -4089 -> -3906.6667
$ws.Cells.Item(77, 8).Value = 884.35895  # H77: 787.13043 -> 884.35895
$ws.Cells.Item(77, 9).Value = 502.06668  # I77: 460 -> 502.06668
$ws.Cells.Item(77, 10).Value = 2158.6667  # J77: 2341 -> 2158.6667
$ws.Cells.Item(77, 11).Value = 2510.3334  # K77: 2300 -> 2510.3334
$ws.Cells.Item(77, 12).Value = 10793.3335  # L77: 11705 -> 10793.3335
$ws.Cells.Item(77, 13).Value = 1857.6666  # M77: 2068 -> 1857.6666
$ws.Cells.Item(77, 14).Value = -19529.3335  # N77: -20441 -> -19529.3335
$ws.Cells.Item(132, 8).Value = 1576.5135  # H132: 1576.7567 -> 1576.5135
$ws.Cells.Item(132, 9).Value = 1265.3214  # I132: 1265.6428 -> 1265.3214
$ws.Cells.Item(132, 11).Value = 3795.9642  # K132: 3796.9284 -> 3795.9642
$ws.Cells.Item(132, 13).Value = -1265.9642  # M132: -1266.9284 -> -1265.9642
$ws.Cells.Item(136, 8).Value = 990  # H136: 945.6667 -> 990
$ws.Cells.Item(136, 9).Value = 990  # I136: 945.6667 -> 990
$ws.Cells.Item(136, 11).Value = 2970  # K136: 2837.0001 -> 2970
$ws.Cells.Item(136, 13).Value = -420  # M136: -287.0001000000002 -> -420

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(20, 8).Value = 2472.1428  # H20: 2331 -> 2472.1428
$ws.Cells.Item(20, 9).Value = 2451.375  # I20: 2218.5 -> 2451.375
$ws.Cells.Item(20, 10).Value = 2499.8333  # J20: 2499.75 -> 2499.8333
$ws.Cells.Item(20, 11).Value = 2451.375  # K20: 2218.5 -> 2451.375
$ws.Cells.Item(20, 12).Value = 2499.8333  # L20: 2499.75 -> 2499.8333
$ws.Cells.Item(20, 13).Value = -2204.375  # M20: -1971.5 -> -2204.375
$ws.Cells.Item(20, 14).Value = -2993.8333  # N20: -2993.75 -> -2993.8333
$ws.Cells.Item(94, 8).Value = 35716052  # H94: 50002004 -> 35716052
$ws.Cells.Item(94, 9).Value = 41668228  # I94: 83334670 -> 41668228
$ws.Cells.Item(94, 10).Value = 3009  # J94: 3009.5 -> 3009
$ws.Cells.Item(94, 11).Value = 41668228  # K94: 83334670 -> 41668228
$ws.Cells.Item(94, 12).Value = 3009  # L94: 3009.5 -> 3009
$ws.Cells.Item(94, 13).Value = -41667777  # M94: -83334219 -> -41667777
$ws.Cells.Item(94, 14).Value = -3911  # N94: -3911.5 -> -3911
$ws.Cells.Item(107, 8).Value = 1050.4783  # H107: 1431.5883 -> 1050.4783
$ws.Cells.Item(107, 9).Value = 733.82355  # I107: 934.25 -> 733.82355
$ws.Cells.Item(107, 10).Value = 1947.6666  # J107: 2625.2 -> 1947.6666
$ws.Cells.Item(107, 11).Value = 733.82355  # K107: 934.25 -> 733.82355
$ws.Cells.Item(107, 12).Value = 1947.6666  # L107: 2625.2 -> 1947.6666
$ws.Cells.Item(107, 13).Value = 1186.17645  # M107: 985.75 -> 1186.17645
$ws.Cells.Item(107, 14).Value = -5787.6666  # N107: -6465.2 -> -5787.6666
$ws.Cells.Item(134, 8).Value = 10916.23  # H134: 6196.3335 -> 10916.23
$ws.Cells.Item(134, 9).Value = 1323.5555  # I134: 935.65 -> 1323.5555
$ws.Cells.Item(134, 11).Value = 3970.6665  # K134: 2806.95 -> 3970.6665
$ws.Cells.Item(134, 13).Value = -1435.6665  # M134: -271.9499999999998 -> -1435.6665

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 1523.2142  # H31: 1291.5 -> 1523.2142
$ws.Cells.Item(31, 9).Value = 1333.381  # I31: 1024.4286 -> 1333.381
$ws.Cells.Item(31, 11).Value = 1333.381  # K31: 1024.4286 -> 1333.381
$ws.Cells.Item(31, 13).Value = -1038.381  # M31: -729.4286 -> -1038.381
$ws.Cells.Item(34, 8).Value = 1523.2142  # H34: 1291.5 -> 1523.2142
$ws.Cells.Item(34, 9).Value = 1333.381  # I34: 1024.4286 -> 1333.381
$ws.Cells.Item(34, 11).Value = 1333.381  # K34: 1024.4286 -> 1333.381
$ws.Cells.Item(34, 13).Value = -1131.381  # M34: -822.4286 -> -1131.381
$ws.Cells.Item(58, 8).Value = 1660.3478  # H58: 1787.65 -> 1660.3478
$ws.Cells.Item(58, 9).Value = 1377.1177  # I58: 1450.2 -> 1377.1177
$ws.Cells.Item(58, 10).Value = 2462.8333  # J58: 2800 -> 2462.8333
$ws.Cells.Item(58, 11).Value = 1377.1177  # K58: 1450.2 -> 1377.1177
$ws.Cells.Item(58, 12).Value = 2462.8333  # L58: 2800 -> 2462.8333
$ws.Cells.Item(58, 13).Value = -1174.1177  # M58: -1247.2 -> -1174.1177
$ws.Cells.Item(58, 14).Value = -2868.8333  # N58: -3206 -> -2868.8333
$ws.Cells.Item(132, 8).Value = 5342.6  # H132: 6845.952 -> 5342.6
$ws.Cells.Item(132, 9).Value = 7139.0557  # I132: 9547.416999999999 -> 7139.0557
$ws.Cells.Item(132, 10).Value = 2647.9167  # J132: 3244 -> 2647.9167
$ws.Cells.Item(132, 11).Value = 21417.1671  # K132: 28642.251 -> 21417.1671
$ws.Cells.Item(132, 12).Value = 7943.750100000001  # L132: 9732 -> 7943.750100000001
$ws.Cells.Item(132, 13).Value = -18887.1671  # M132: -26112.251 -> -18887.1671
$ws.Cells.Item(132, 14).Value = -13003.7501  # N132: -14792 -> -13003.7501
$ws.Cells.Item(134, 8).Value = 1558.5278  # H134: 1665.6061 -> 1558.5278
$ws.Cells.Item(134, 9).Value = 1538.9667  # I134: 1667.6666 -> 1538.9667
$ws.Cells.Item(134, 11).Value = 4616.9001  # K134: 5002.9998 -> 4616.9001
$ws.Cells.Item(134, 13).Value = -2081.9001  # M134: -2467.9998 -> -2081.9001
$ws.Cells.Item(136, 8).Value = 1660.3478  # H136: 1787.65 -> 1660.3478
$ws.Cells.Item(136, 9).Value = 1377.1177  # I136: 1450.2 -> 1377.1177
$ws.Cells.Item(136, 10).Value = 2462.8333  # J136: 2800 -> 2462.8333
$ws.Cells.Item(136, 11).Value = 4131.3531  # K136: 4350.6 -> 4131.3531
$ws.Cells.Item(136, 12).Value = 7388.499899999999  # L136: 8400 -> 7388.499899999999
$ws.Cells.Item(136, 13).Value = -1581.3531  # M136: -1800.6 -> -1581.3531
$ws.Cells.Item(136, 14).Value = -12488.4999  # N136: -13500 -> -12488.4999

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(5, 8).Value = 601.37933  # H5: 604.65515 -> 601.37933
$ws.Cells.Item(5, 9).Value = 593.75  # I5: 603.5217 -> 593.75
$ws.Cells.Item(5, 10).Value = 638  # J5: 609 -> 638
$ws.Cells.Item(5, 11).Value = 1781.25  # K5: 1810.5651 -> 1781.25
$ws.Cells.Item(5, 12).Value = 1914  # L5: 1827 -> 1914
$ws.Cells.Item(5, 13).Value = -1669.25  # M5: -1698.5651 -> -1669.25
$ws.Cells.Item(5, 14).Value = -2138  # N5: -2051 -> -2138
$ws.Cells.Item(9, 8).Value = 500  # H9: 1196 -> 500
$ws.Cells.Item(9, 10).Value = 500  # J9: 1196 -> 500
$ws.Cells.Item(9, 12).Value = 1500  # L9: 3588 -> 1500
$ws.Cells.Item(9, 14).Value = -1948  # N9: -4036 -> -1948
$ws.Cells.Item(39, 8).Value = 1796.9584  # H39: 1899.4546 -> 1796.9584
$ws.Cells.Item(39, 10).Value = 1510.8096  # J39: 1599.3684 -> 1510.8096
$ws.Cells.Item(39, 12).Value = 4532.4288  # L39: 4798.1052 -> 4532.4288
$ws.Cells.Item(39, 14).Value = -5120.4288  # N39: -5386.1052 -> -5120.4288
$ws.Cells.Item(115, 8).Value = 4256.75  # H115: 4342.6665 -> 4256.75
$ws.Cells.Item(115, 10).Value = 4999.6665  # J115: 5500 -> 4999.6665
$ws.Cells.Item(115, 12).Value = 14998.9995  # L115: 16500 -> 14998.9995
$ws.Cells.Item(115, 14).Value = -17348.9995  # N115: -18850 -> -17348.9995
$ws.Cells.Item(121, 8).Value = 740  # H121: 762.4666999999999 -> 740
$ws.Cells.Item(121, 9).Value = 458.8  # I121: 473.5 -> 458.8
$ws.Cells.Item(121, 10).Value = 867.8182  # J121: 867.5454999999999 -> 867.8182
$ws.Cells.Item(121, 11).Value = 1376.4  # K121: 1420.5 -> 1376.4
$ws.Cells.Item(121, 12).Value = 2603.4546  # L121: 2602.6365 -> 2603.4546
$ws.Cells.Item(121, 13).Value = -66.40000000000009  # M121: -110.5 -> -66.40000000000009
$ws.Cells.Item(121, 14).Value = -5223.4546  # N121: -5222.6365 -> -5223.4546
$ws.Cells.Item(122, 8).Value = 815.2857  # H122: 741.7778 -> 815.2857
$ws.Cells.Item(122, 9).Value = 665.6667  # I122: 517.2 -> 665.6667
$ws.Cells.Item(122, 10).Value = 927.5  # J122: 1022.5 -> 927.5
$ws.Cells.Item(122, 11).Value = 5991.0003  # K122: 4654.8 -> 5991.0003
$ws.Cells.Item(122, 12).Value = 8347.5  # L122: 9202.5 -> 8347.5
$ws.Cells.Item(122, 13).Value = -3541.0003  # M122: -2204.8 -> -3541.0003
$ws.Cells.Item(122, 14).Value = -13247.5  # N122: -14102.5 -> -13247.5
$ws.Cells.Item(132, 8).Value = 1316  # H132: 1998.25 -> 1316
$ws.Cells.Item(132, 9).Value = 633.3333  # I132: 888 -> 633.3333
$ws.Cells.Item(132, 10).Value = 1608.5714  # J132: 2368.3333 -> 1608.5714
$ws.Cells.Item(132, 11).Value = 5699.9997  # K132: 7992 -> 5699.9997
$ws.Cells.Item(132, 12).Value = 14477.1426  # L132: 21314.9997 -> 14477.1426
$ws.Cells.Item(132, 13).Value = -3169.9997  # M132: -5462 -> -3169.9997
$ws.Cells.Item(132, 14).Value = -19537.1426  # N132: -26374.9997 -> -19537.1426
$ws.Cells.Item(135, 8).Value = 601.37933  # H135: 604.65515 -> 601.37933
$ws.Cells.Item(135, 9).Value = 593.75  # I135: 603.5217 -> 593.75
$ws.Cells.Item(135, 10).Value = 638  # J135: 609 -> 638
$ws.Cells.Item(135, 11).Value = 5343.75  # K135: 5431.6953 -> 5343.75
$ws.Cells.Item(135, 12).Value = 5742  # L135: 5481 -> 5742
$ws.Cells.Item(135, 13).Value = -2808.75  # M135: -2896.6953 -> -2808.75
$ws.Cells.Item(135, 14).Value = -10812  # N135: -10551 -> -10812

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(15, 8).Value = 0  # H15: 20000 -> 0
$ws.Cells.Item(15, 10).Value = 0  # J15: 20000 -> 0
$ws.Cells.Item(15, 12).Value = 0  # L15: 20000 -> 0
$ws.Cells.Item(15, 14).Value = ""  # N15: clear (was -20576)
$ws.Cells.Item(70, 8).Value = 56252876  # H70: 56253004 -> 56252876
$ws.Cells.Item(70, 9).Value = 62502750  # I70: 125002000 -> 62502750
$ws.Cells.Item(70, 10).Value = 50003004  # J70: 33336672 -> 50003004
$ws.Cells.Item(70, 11).Value = 62502750  # K70: 125002000 -> 62502750
$ws.Cells.Item(70, 12).Value = 50003004  # L70: 33336672 -> 50003004
$ws.Cells.Item(70, 13).Value = -62502480  # M70: -125001730 -> -62502480
$ws.Cells.Item(70, 14).Value = -50003544  # N70: -33337212 -> -50003544
$ws.Cells.Item(73, 8).Value = 56252876  # H73: 56253004 -> 56252876
$ws.Cells.Item(73, 9).Value = 62502750  # I73: 125002000 -> 62502750
$ws.Cells.Item(73, 10).Value = 50003004  # J73: 33336672 -> 50003004
$ws.Cells.Item(73, 11).Value = 62502750  # K73: 125002000 -> 62502750
$ws.Cells.Item(73, 12).Value = 50003004  # L73: 33336672 -> 50003004
$ws.Cells.Item(73, 13).Value = -62501814  # M73: -125001064 -> -62501814
$ws.Cells.Item(73, 14).Value = -50004876  # N73: -33338544 -> -50004876
$ws.Cells.Item(81, 8).Value = 0  # H81: 20000 -> 0
$ws.Cells.Item(81, 10).Value = 0  # J81: 20000 -> 0
$ws.Cells.Item(81, 12).Value = 0  # L81: 20000 -> 0
$ws.Cells.Item(81, 14).Value = ""  # N81: clear (was -21996)
$ws.Cells.Item(84, 8).Value = 0  # H84: 20000 -> 0
$ws.Cells.Item(84, 10).Value = 0  # J84: 20000 -> 0
$ws.Cells.Item(84, 12).Value = 0  # L84: 60000 -> 0
$ws.Cells.Item(84, 14).Value = ""  # N84: clear (was -69984)
$ws.Cells.Item(122, 8).Value = 1385.75  # H122: 1482.6666 -> 1385.75
$ws.Cells.Item(122, 9).Value = 1488.4667  # I122: 1638.4615 -> 1488.4667
$ws.Cells.Item(122, 11).Value = 4465.4001  # K122: 4915.3845 -> 4465.4001
$ws.Cells.Item(122, 13).Value = -2015.4001  # M122: -2465.3845 -> -2015.4001
$ws.Cells.Item(132, 8).Value = 2021.9032  # H132: 2173.6667 -> 2021.9032
$ws.Cells.Item(132, 9).Value = 1712.8096  # I132: 1881.1177 -> 1712.8096
$ws.Cells.Item(132, 11).Value = 5138.4288  # K132: 5643.3531 -> 5138.4288
$ws.Cells.Item(132, 13).Value = -2608.4288  # M132: -3113.3531 -> -2608.4288

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(8, 8).Value = 0  # H8: 25000 -> 0
$ws.Cells.Item(8, 10).Value = 0  # J8: 25000 -> 0
$ws.Cells.Item(8, 12).Value = 0  # L8: 25000 -> 0
$ws.Cells.Item(8, 14).Value = ""  # N8: clear (was -25280)
$ws.Cells.Item(112, 8).Value = 45999.445  # H112: 48499 -> 45999.445
$ws.Cells.Item(112, 10).Value = 45999.445  # J112: 48499 -> 45999.445
$ws.Cells.Item(112, 12).Value = 45999.445  # L112: 48499 -> 45999.445
$ws.Cells.Item(112, 14).Value = -48953.445  # N112: -51453 -> -48953.445
$ws.Cells.Item(132, 8).Value = 33362.594  # H132: 34410.16 -> 33362.594
$ws.Cells.Item(132, 9).Value = 1749.3572  # I132: 1815.6154 -> 1749.3572
$ws.Cells.Item(132, 11).Value = 5248.071599999999  # K132: 5446.8462 -> 5248.071599999999
$ws.Cells.Item(132, 13).Value = -2718.071599999999  # M132: -2916.8462 -> -2718.071599999999
$ws.Cells.Item(136, 8).Value = 5432.773  # H136: 6221.2104 -> 5432.773
$ws.Cells.Item(136, 9).Value = 7501.1333  # I136: 9266.583000000001 -> 7501.1333
$ws.Cells.Item(136, 11).Value = 22503.3999  # K136: 27799.749 -> 22503.3999
$ws.Cells.Item(136, 13).Value = -19953.3999  # M136: -25249.749 -> -19953.3999
$ws.Cells.Item(140, 8).Value = 43856  # H140: 43868 -> 43856
$ws.Cells.Item(140, 10).Value = 43856  # J140: 43868 -> 43856
$ws.Cells.Item(140, 12).Value = 43856  # L140: 43868 -> 43856
$ws.Cells.Item(140, 14).Value = -54216  # N140: -54228 -> -54216

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(92, 8).Value = 28058.166  # H92: 29725 -> 28058.166
$ws.Cells.Item(92, 10).Value = 28058.166  # J92: 29725 -> 28058.166
$ws.Cells.Item(92, 12).Value = 28058.166  # L92: 29725 -> 28058.166
$ws.Cells.Item(92, 14).Value = -33050.166  # N92: -34717 -> -33050.166
$ws.Cells.Item(132, 8).Value = 6495.909  # H132: 6745.5 -> 6495.909
$ws.Cells.Item(132, 9).Value = 6651.143  # I132: 7093 -> 6651.143
$ws.Cells.Item(132, 11).Value = 19953.429  # K132: 21279 -> 19953.429
$ws.Cells.Item(132, 13).Value = -17423.429  # M132: -18749 -> -17423.429
$ws.Cells.Item(136, 8).Value = 735.2105  # H136: 721.6667 -> 735.2105
$ws.Cells.Item(136, 10).Value = 1011.8  # J136: 1020 -> 1011.8
$ws.Cells.Item(136, 12).Value = 3035.4  # L136: 3060 -> 3035.4
$ws.Cells.Item(136, 14).Value = -8135.4  # N136: -8160 -> -8135.4

